$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextPlaceholder($ref, $text) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $ws.Range("A14").Copy()
    $cell.PasteSpecial(-4122)
}

function Set-NumericStyle($ref, $value) {
    $cell = $ws.Range($ref)
    $cell.Value = $value
    $ws.Range("G15").Copy()
    $cell.PasteSpecial(-4122)
}

# --- Header text updates (Volume number, report week dates) ---
$ws.Range("A8").Value = "Volume 29   Number  47"
$ws.Range("C9").Value = "Report Covering the Week  11/21/2022  Through  11/27/2022"

# --- Cells that flip from numeric to text placeholder ("0" / "***.*") ---
Set-TextPlaceholder "C15" "0"
Set-TextPlaceholder "D22" "0"
Set-TextPlaceholder "E22" "***.*"
Set-TextPlaceholder "C26" "0"
Set-TextPlaceholder "D27" "0"
Set-TextPlaceholder "E27" "***.*"
Set-TextPlaceholder "D28" "0"
Set-TextPlaceholder "E28" "***.*"
Set-TextPlaceholder "D29" "0"
Set-TextPlaceholder "E29" "***.*"

# --- Cells that flip from text placeholder to numeric ---
Set-NumericStyle "C18" 2
Set-NumericStyle "C27" 2

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("N14").Value = -82.352941176470
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = -100
$ws.Range("F15").Value = 1
$ws.Range("H15").Value = -83.333333333333
$ws.Range("J15").Value = 31
$ws.Range("K15").Value = 16.129032258064
$ws.Range("N15").Value = 12.5
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 37
$ws.Range("H16").Value = 68.181818181818
$ws.Range("I16").Value = 287
$ws.Range("J16").Value = 222
$ws.Range("K16").Value = 29.279279279279
$ws.Range("L16").Value = 47.938144329896
$ws.Range("M16").Value = -10.591900311526
$ws.Range("N16").Value = -75.943000838223
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 26
$ws.Range("G17").Value = 30
$ws.Range("H17").Value = -13.333333333333
$ws.Range("I17").Value = 382
$ws.Range("J17").Value = 353
$ws.Range("K17").Value = 8.215297450424
$ws.Range("L17").Value = 35.943060498220
$ws.Range("M17").Value = 28.187919463087
$ws.Range("N17").Value = -2.551020408163
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -9.090909090909
$ws.Range("I18").Value = 135
$ws.Range("J18").Value = 128
$ws.Range("K18").Value = 5.46875
$ws.Range("L18").Value = -10.596026490066
$ws.Range("M18").Value = -49.814126394052
$ws.Range("N18").Value = -92.5
$ws.Range("C19").Value = 17
$ws.Range("D19").Value = 66
$ws.Range("E19").Value = -74.242424242424
$ws.Range("F19").Value = 74
$ws.Range("G19").Value = 123
$ws.Range("H19").Value = -39.837398373983
$ws.Range("I19").Value = 891
$ws.Range("J19").Value = 577
$ws.Range("K19").Value = 54.419410745234
$ws.Range("L19").Value = 84.472049689441
$ws.Range("M19").Value = 94.967177242888
$ws.Range("N19").Value = -32.5
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 133.333333333333
$ws.Range("F20").Value = 28
$ws.Range("G20").Value = 16
$ws.Range("H20").Value = 75
$ws.Range("I20").Value = 290
$ws.Range("J20").Value = 190
$ws.Range("K20").Value = 52.631578947368
$ws.Range("L20").Value = 54.255319148936
$ws.Range("M20").Value = 34.259259259259
$ws.Range("N20").Value = -85.963213939980
$ws.Range("C21").Value = 38
$ws.Range("D21").Value = 89
$ws.Range("E21").Value = -57.303370786516
$ws.Range("F21").Value = 176
$ws.Range("G21").Value = 208
$ws.Range("H21").Value = -15.384615384615
$ws.Range("I21").Value = 2024
$ws.Range("J21").Value = 1505
$ws.Range("K21").Value = 34.485049833887
$ws.Range("L21").Value = 51.610486891385
$ws.Range("M21").Value = 27.215587680704
$ws.Range("N21").Value = -70.322580645161
$ws.Range("F22").Value = 4
$ws.Range("H22").Value = 33.333333333333
$ws.Range("C24").Value = 36
$ws.Range("D24").Value = 35
$ws.Range("E24").Value = 2.857142857142
$ws.Range("F24").Value = 187
$ws.Range("G24").Value = 126
$ws.Range("H24").Value = 48.412698412698
$ws.Range("I24").Value = 1765
$ws.Range("J24").Value = 1253
$ws.Range("K24").Value = 40.861931364724
$ws.Range("L24").Value = 37.998436278342
$ws.Range("M24").Value = 71.359223300970
$ws.Range("C25").Value = 22
$ws.Range("D25").Value = 19
$ws.Range("E25").Value = 15.789473684210
$ws.Range("F25").Value = 68
$ws.Range("G25").Value = 74
$ws.Range("H25").Value = -8.108108108108
$ws.Range("I25").Value = 824
$ws.Range("J25").Value = 779
$ws.Range("K25").Value = 5.776636713735
$ws.Range("L25").Value = 19.767441860465
$ws.Range("M25").Value = -2.830188679245
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = 2
$ws.Range("J26").Value = 52
$ws.Range("K26").Value = -7.692307692307
$ws.Range("L26").Value = 11.627906976744
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 260
$ws.Range("I27").Value = 104
$ws.Range("K27").Value = 8.333333333333
$ws.Range("L27").Value = 38.666666666666
$ws.Range("N28").Value = -85.454545454545
$ws.Range("N29").Value = -86.274509803921
